$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187 (pushes the existing row 187..215 down to 188..216)
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new weekly record
$ws.Range("A187").Value = 10
$ws.Range("B187").Value = "Vega Modelo de Temuco"
$ws.Range("C187").Value = "La Araucanía"
$ws.Range("D187").Value = 44474
$ws.Range("D187").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E187").Value = 9
$ws.Range("F187").Value = 100112037
$ws.Range("G187").Value = "Cebollín"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 90
$ws.Range("K187").Value = 5000
$ws.Range("L187").Value = 8000
$ws.Range("M187").Value = 6444
$ws.Range("N187").Value = "`$/docena de paquetes"
$ws.Range("O187").Value = "Provincia de Cautín"
$ws.Range("P187").Value = 537
$ws.Range("Q187").Value = 12
$ws.Range("R187").Value = "Hortaliza"
